$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.448.88'
$ws.Range('E2').Value = '  +0.51%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.329.03'
$ws.Range('E3').Value = '  -0.71%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.29'
$ws.Range('E5').Value = '  -1.23%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.31'
$ws.Range('E6').Value = '  -2.30%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.509'
$ws.Range('E7').Value = '  -2.74%  '

# Row 8
$ws.Range('E8').Value = '  +0.15%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').Value = '  -2.33%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.31'
$ws.Range('E10').Value = '  -1.65%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').Value = '  -1.50%  '

# Row 12
$ws.Range('E12').Value = '  +0.45%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.78'
$ws.Range('E13').Value = '  -2.62%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.692.55'
$ws.Range('E14').Value = '  -0.40%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.63'
$ws.Range('E15').Value = '  -0.26%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.328.90'
$ws.Range('E16').Value = '  -0.75%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  -0.28%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.375.83'
$ws.Range('E18').Value = '  +0.40%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.92'
$ws.Range('E19').Value = '  +0.11%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0910'
$ws.Range('E20').Value = '  -1.54%  '

# Row 21
$ws.Range('E21').Value = '  -2.48%  '

# Row 22
$ws.Range('E22').Value = '  +0.24%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.01'
$ws.Range('E23').Value = '  -1.37%  '

# Row 24
$ws.Range('E24').Value = '  -3.35%  '

# Row 25
$ws.Range('E25').Value = '  -2.59%  '

# Row 26
$ws.Range('E26').Value = '  -0.01%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.96'
$ws.Range('E27').Value = '  -1.41%  '

# Row 28
$ws.Range('E28').Value = '  -1.12%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.58'
$ws.Range('E29').Value = '  -4.86%  '

# Row 30
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '164.97'
$ws.Range('E30').Value = '  +1.54%  '

# Row 31
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.20'
$ws.Range('E31').Value = '  -3.79%  '

# Row 32
$ws.Range('E32').Value = '  +0.09%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.05'
$ws.Range('E33').Value = '  -3.95%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.57'
$ws.Range('E34').Value = '  +0.03%  '

# Row 35
$ws.Range('E35').Value = '  -4.88%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.06'
$ws.Range('E36').Value = '  -6.13%  '

# Row 37
$ws.Range('E37').Value = '  -4.19%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.91'
$ws.Range('E38').Value = '  -5.79%  '

# Row 39
$ws.Range('E39').Value = '  -4.07%  '

# Row 40
$ws.Range('E40').Value = '  -3.86%  '

# Row 41
$ws.Range('E41').Value = '  -2.97%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.42'
$ws.Range('E42').Value = '  -0.69%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.978.34'
$ws.Range('E43').Value = '  -0.25%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0284'
$ws.Range('E44').Value = '  -2.29%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.43'
$ws.Range('E45').Value = '  -5.83%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.08'
$ws.Range('E46').Value = '  -3.08%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('E47').Value = '  -4.25%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.80'
$ws.Range('E48').Value = '  -4.63%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.80'
$ws.Range('E49').Value = '  +2.54%  '

# Row 50
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.553.46'
$ws.Range('E50').Value = '  +0.22%  '

# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.56'
$ws.Range('E51').Value = '  -1.85%  '
